$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: title row (merged A1:D1). Keep "标题" text, re-apply centered
# alignment + font so Excel regenerates the font/style table the same way
# the real edit did (duplicate default font becomes fontId 2).
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").Font.Name = "等线"
$ws.Range("A1:D1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:D1").VerticalAlignment = -4108     # xlCenter

$ws.Range("F1:H1").Font.Name = "等线"
$ws.Range("F1:H1").HorizontalAlignment = -4108
$ws.Range("F1:H1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 2: class / department line.
# A2 = 班级 (label), B2 = 初三十班 (class name)
# C2 = 学部 (label), D2 = 初中部 (department)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "班级"
$ws.Range("B2").Value = "初三十班"
$ws.Range("C2").Value = "学部"
$ws.Range("D2").Value = "初中部"

$ws.Range("A2:H2").Font.Name = "等线"
$ws.Range("A2:H2").HorizontalAlignment = -4108
$ws.Range("A2:H2").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 3: column headers for the roster table -> 姓名 (name) / 方向 (hometown)
# Clear the old C3/D3 (大小周 / 方向 leftovers) since the new layout only
# uses columns A and B.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "姓名"
$ws.Range("B3").Value = "方向"
$ws.Range("C3:D3").Clear() | Out-Null

$ws.Range("A3:B3").Font.Name = "等线"
$ws.Range("A3:B3").VerticalAlignment = -4108
$ws.Range("A3:B3").HorizontalAlignment = -4142   # xlGeneral (no forced centering)

# ---------------------------------------------------------------------------
# Roster data rows 4-13: 姓名 (name) / 方向 (hometown)
# ---------------------------------------------------------------------------
$names = @("刘浩宇","刘壮","李沚璠","谭一诺","韩蓝颖","李月童","崔兆言","万宇","李大千","张鹏")
$towns = @("沧州","廊坊","沧州","任丘","沧州","沧州","沧州","沧州","张家口","张家口")

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $towns[$i]
}

# Clear the old leftover values in columns C/D for rows 4 and 5 (originally
# held A方向/true and A方/TRE, no longer part of the new layout).
$ws.Range("C4:D5").Clear() | Out-Null

$ws.Range("A4:B13").Font.Name = "等线"
$ws.Range("A4:B13").VerticalAlignment = -4108
$ws.Range("A4:B13").HorizontalAlignment = -4142

for ($r = 4; $r -le 13; $r++) {
    $ws.Rows.Item($r).RowHeight = 18
}

# ---------------------------------------------------------------------------
# Selection matches the author's saved cursor position (A14) after the edit.
# ---------------------------------------------------------------------------
$ws.Range("A14").Select() | Out-Null
